# Add a new worksheet "ODI Batting Extra" right after "ODI Batting",
# matching the header style of the existing sheets, and populate it
# with per-innings extra-stats data keyed by MATCH_CODE.

$wb = $excel.ActiveWorkbook
$wsBatting = $wb.Worksheets.Item(2)

# Insert the new sheet directly after the "ODI Batting" sheet.
$ws = $wb.Worksheets.Add($null, $wsBatting)
$ws.Name = "ODI Batting Extra"

# --- Header row -------------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Match the bold/bordered/centered header formatting used on the other sheets
# by copying it from the "ODI Batting" header row.
$wsBatting.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Data rows ----------------------------------------------------------
# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$data = @(
    @("3860", "", "", "", "", "NO"),
    @("3862", "", "", "", "", "NO"),
    @("4138", "", "", "", "", "NO"),
    @("4139", "", "", "", "", "NO"),
    @("4149", 5,  "0", "0", "", "NO"),
    @("4406", "", "", "", "", "NO"),
    @("4625", 4,  "6", "7", "32.90%", "YES")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    # MATCH_CODE: numeric-looking text -> force text storage with a
    # leading apostrophe so it round-trips as a string, not a number.
    $ws.Cells.Item($row, 1).Formula = "'" + $rec[0]

    # BATTING_POSITION: a genuine number when present, otherwise blank.
    if ($rec[1] -ne "") {
        $ws.Cells.Item($row, 2).Value = $rec[1]
    }

    # NUM_4 / NUM_6: numeric-looking text -> force text storage.
    if ($rec[2] -ne "") {
        $ws.Cells.Item($row, 3).Formula = "'" + $rec[2]
    }
    if ($rec[3] -ne "") {
        $ws.Cells.Item($row, 4).Formula = "'" + $rec[3]
    }

    # PERCENT_RUNS_OF_TOTAL: literal text (e.g. "32.90%"), not a percentage number.
    if ($rec[4] -ne "") {
        $ws.Cells.Item($row, 5).Formula = "'" + $rec[4]
    }

    # MAN_OF_MATCH: plain YES/NO text.
    $ws.Cells.Item($row, 6).Value = $rec[5]
}
